$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 2.3
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 2.88
$ws.Range("K4").Value = 2.3
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 13
$ws.Range("Q4").Value = 1.7
$ws.Range("R4").Value = 2.1
$ws.Range("Z4").Value = 21
$ws.Range("AO4").Value = 12
$ws.Range("AP4").Value = 19
$ws.Range("BA4").Value = 67

# Row 5
$ws.Range("G5").Value = 2.1
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 3.5
$ws.Range("J5").Value = 2.75
$ws.Range("L5").Value = 3.75
$ws.Range("N5").Value = 13
$ws.Range("S5").Value = 1.36
$ws.Range("T5").Value = 3
$ws.Range("Z5").Value = 19
$ws.Range("AI5").Value = 19
$ws.Range("AL5").Value = 26
$ws.Range("AM5").Value = 29
$ws.Range("AQ5").Value = 41
$ws.Range("AT5").Value = 3
$ws.Range("AX5").Value = 17

# Row 6
$ws.Range("G6").Value = 2.7
$ws.Range("H6").Value = 3.25
$ws.Range("I6").Value = 2.55
$ws.Range("J6").Value = 3.4
$ws.Range("K6").Value = 2.1
$ws.Range("L6").Value = 3.25
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 10
$ws.Range("O6").Value = 1.33
$ws.Range("P6").Value = 3.25
$ws.Range("Q6").Value = 2.05
$ws.Range("R6").Value = 1.75
$ws.Range("S6").Value = 1.44
$ws.Range("T6").Value = 2.63
$ws.Range("U6").Value = 1.8
$ws.Range("V6").Value = 1.91
$ws.Range("W6").Value = 8.5
$ws.Range("X6").Value = 13
$ws.Range("AB6").Value = 34
$ws.Range("AC6").Value = 9.5
$ws.Range("AD6").Value = 6.5
$ws.Range("AE6").Value = 15
$ws.Range("AF6").Value = 51
$ws.Range("AG6").Value = 251
$ws.Range("AH6").Value = 8.5
$ws.Range("AI6").Value = 12
$ws.Range("AK6").Value = 26
$ws.Range("AL6").Value = 21
$ws.Range("AM6").Value = 29
$ws.Range("AP6").Value = 26
$ws.Range("AQ6").Value = 51
$ws.Range("AR6").Value = 67
$ws.Range("AS6").Value = 201
$ws.Range("AT6").Value = 2.63
$ws.Range("AU6").Value = 8
$ws.Range("AW6").Value = 4.5
$ws.Range("AX6").Value = 15
$ws.Range("AY6").Value = 23
$ws.Range("AZ6").Value = 51
$ws.Range("BA6").Value = 67
$ws.Range("BB6").Value = 151

# Row 8
$ws.Range("G8").Value = 1.93
$ws.Range("H8").Value = 3.35
$ws.Range("I8").Value = 3.65
$ws.Range("J8").Value = 2.47
$ws.Range("K8").Value = 2.15
$ws.Range("L8").Value = 4.05
$ws.Range("M8").Value = 1.06
$ws.Range("N8").Value = 8.42
$ws.Range("O8").Value = 1.34
$ws.Range("P8").Value = 2.72
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 1.65
$ws.Range("S8").Value = 1.42
$ws.Range("T8").Value = 2.47
$ws.Range("W8").Value = 6.4
$ws.Range("X8").Value = 8.5
$ws.Range("Z8").Value = 16.5
$ws.Range("AA8").Value = 16.5
$ws.Range("AC8").Value = 8.75
$ws.Range("AD8").Value = 6.5
$ws.Range("AE8").Value = 16.5
$ws.Range("AH8").Value = 9.5
$ws.Range("AI8").Value = 18.5
$ws.Range("AJ8").Value = 13
$ws.Range("AK8").Value = 55
$ws.Range("AL8").Value = 37
$ws.Range("AN8").Value = 3.75
$ws.Range("AO8").Value = 9.5
$ws.Range("AP8").Value = 18.5
$ws.Range("AQ8").Value = 35
$ws.Range("AR8").Value = 65
$ws.Range("AT8").Value = 2.6
$ws.Range("AU8").Value = 7.2
$ws.Range("AV8").Value = 65
$ws.Range("AW8").Value = 5.4
$ws.Range("AX8").Value = 19.5
$ws.Range("AY8").Value = 26

# Row 10
$ws.Range("G10").Value = 5.8
$ws.Range("H10").Value = 4.35
$ws.Range("I10").Value = 1.42
$ws.Range("J10").Value = 5.5
$ws.Range("K10").Value = 2.4
$ws.Range("L10").Value = 1.88
$ws.Range("N10").Value = 12.7
$ws.Range("P10").Value = 4.05
$ws.Range("Q10").Value = 1.57
$ws.Range("R10").Value = 2.12
$ws.Range("S10").Value = 1.3
$ws.Range("T10").Value = 3.32
$ws.Range("U10").Value = 1.76
$ws.Range("V10").Value = 2.01
$ws.Range("W10").Value = 15
$ws.Range("X10").Value = 30
$ws.Range("Y10").Value = 15.5
$ws.Range("Z10").Value = 80
$ws.Range("AA10").Value = 45
$ws.Range("AB10").Value = 40
$ws.Range("AD10").Value = 7.6
$ws.Range("AE10").Value = 14.5
$ws.Range("AF10").Value = 55
$ws.Range("AG10").Value = 350
$ws.Range("AH10").Value = 6.8
$ws.Range("AI10").Value = 6.3
$ws.Range("AJ10").Value = 7.1
$ws.Range("AK10").Value = 8.25
$ws.Range("AL10").Value = 9.25
$ws.Range("AM10").Value = 18.5
$ws.Range("AN10").Value = 7.4
$ws.Range("AO10").Value = 32
$ws.Range("AP10").Value = 35
$ws.Range("AQ10").Value = 200
$ws.Range("AR10").Value = 200
$ws.Range("AS10").Value = 400
$ws.Range("AT10").Value = 3.15
$ws.Range("AU10").Value = 7.9
$ws.Range("AV10").Value = 65
$ws.Range("AW10").Value = 3.35
$ws.Range("AX10").Value = 6.4
$ws.Range("AY10").Value = 15.5
$ws.Range("AZ10").Value = 18
$ws.Range("BB10").Value = 200

# Row 11
$ws.Range("G11").Value = 32
$ws.Range("J11").Value = 21
$ws.Range("K11").Value = 3.45
$ws.Range("L11").Value = 1.26
$ws.Range("N11").Value = 16.5
$ws.Range("O11").Value = 1.04
$ws.Range("P11").Value = 7
$ws.Range("Q11").Value = 1.23
$ws.Range("R11").Value = 3.42
$ws.Range("S11").Value = 1.14
$ws.Range("T11").Value = 4.9
$ws.Range("U11").Value = 2.45
$ws.Range("V11").Value = 1.52
$ws.Range("X11").Value = 500
$ws.Range("Y11").Value = 120
$ws.Range("AB11").Value = 400
$ws.Range("AC11").Value = 19
$ws.Range("AE11").Value = 45
$ws.Range("AF11").Value = 200
$ws.Range("AH11").Value = 9
$ws.Range("AI11").Value = 5.7
$ws.Range("AJ11").Value = 12
$ws.Range("AK11").Value = 5.1
$ws.Range("AL11").Value = 10.5
$ws.Range("AM11").Value = 40
$ws.Range("AN11").Value = 30
$ws.Range("AO11").Value = 300
$ws.Range("AP11").Value = 150
$ws.Range("AT11").Value = 4.85
$ws.Range("AU11").Value = 13.5
$ws.Range("AV11").Value = 120
$ws.Range("AW11").Value = 3.2
$ws.Range("AX11").Value = 3.8
$ws.Range("AY11").Value = 15.5
$ws.Range("AZ11").Value = 6.5
$ws.Range("BA11").Value = 27
$ws.Range("BB11").Value = 200

# Row 12
$ws.Range("L12").Value = 2.1
$ws.Range("O12").Value = 1.29
$ws.Range("P12").Value = 3.5
$ws.Range("Q12").Value = 1.98
$ws.Range("R12").Value = 1.88
$ws.Range("AB12").Value = 51
$ws.Range("AS12").Value = 500

# Row 13
$ws.Range("M13").Value = 1.05
$ws.Range("N13").Value = 11
$ws.Range("O13").Value = 1.29
$ws.Range("P13").Value = 3.5
$ws.Range("Q13").Value = 2
$ws.Range("R13").Value = 1.85

# Row 15
$ws.Range("Q15").Value = 2.25
$ws.Range("R15").Value = 1.62

# Row 16
$ws.Range("G16").Value = 2.25
$ws.Range("H16").Value = 2.9
$ws.Range("I16").Value = 3.5
$ws.Range("S16").Value = 1.62
$ws.Range("T16").Value = 2.2
$ws.Range("X16").Value = 9.5
$ws.Range("AA16").Value = 23
$ws.Range("AC16").Value = 6
$ws.Range("AH16").Value = 7.5
$ws.Range("AM16").Value = 51
$ws.Range("AO16").Value = 15
$ws.Range("AT16").Value = 2.2
$ws.Range("AW16").Value = 5
$ws.Range("AY16").Value = 41

# Row 17
$ws.Range("G17").Value = 1.9
$ws.Range("I17").Value = 3.3
$ws.Range("L17").Value = 3.5
$ws.Range("P17").Value = 5.4
$ws.Range("Q17").Value = 1.38
$ws.Range("R17").Value = 2.82
$ws.Range("S17").Value = 1.22
$ws.Range("U17").Value = 1.38
$ws.Range("V17").Value = 2.8
$ws.Range("Z17").Value = 19.5
$ws.Range("AB17").Value = 16
$ws.Range("AD17").Value = 9
$ws.Range("AF17").Value = 29
$ws.Range("AG17").Value = 120
$ws.Range("AH17").Value = 19
$ws.Range("AL17").Value = 23
$ws.Range("AR17").Value = 40
$ws.Range("AX17").Value = 16
$ws.Range("BA17").Value = 65
